$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# 1. Insert a new column C ("MO TA" / description) between TEN (B) and GHI CHU (old C, becomes D)
$ws.Columns.Item(3).Insert()

# 2. New column header
$ws.Range("C2").Value = "MÔ TẢ"

# 3. Fill in the MÔ TẢ column for the existing technique rows (3-8)
$ws.Range("C3").Value = "5 mẫu thiết kế"
$ws.Range("C4").Value = "11 mẫu thiết kế"
$ws.Range("C5").Value = "7 mẫu thiết kế"
$ws.Range("C6").Value = "Mô hình Model-Controller-View"
$ws.Range("C7").Value = "Tạo đối tượng bên client khi dùng web serivce"
$ws.Range("C8").Value = "Bất đồng bộ khi gọi các phương thức của web service"

# 4. Update the note for MVC (row 6) in the now-shifted GHI CHU column (D)
$ws.Range("D6").Value = "Khó áp dụng, có thể dùng với PureMVC hoặc MVC#"

# 5. Update technique name + description for row 9 (Delegate & Event -> Delegate & Event (Public & Subscriber))
$ws.Range("B9").Value = "Delegate & Event (Public & Subscriber)"
$ws.Range("C9").Value = "Phần xử lý sự kiện nằm riêng trong một class"

# 6. Insert a new row for "Message Transformation" before the current blank row (row 10),
#    cloning formatting from row 9 first, then overwriting the values that differ.
$ws.Rows.Item(10).Insert()
$ws.Range("A9:D9").Copy($ws.Range("A10:D10"))
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Message Transformation"
$ws.Range("C10").Value = "Đọc file XML và chuyển thành giao diện"
$ws.Range("D10").Value = "Có thể áp dụng"

# 7. Add one more blank row at the bottom (STT 10), cloning formatting from the existing blank row (row 11)
$ws.Rows.Item(12).Insert()
$ws.Range("A11:D11").Copy($ws.Range("A12:D12"))
$ws.Range("A12").Value = 10
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()

# Renumber the old blank row (now row 11) as STT 9
$ws.Range("A11").Value = 9

# Column widths (bestFit-like sizing to roughly match the new content)
$ws.Columns.Item(2).ColumnWidth = 32.8
$ws.Columns.Item(3).ColumnWidth = 44.93
$ws.Columns.Item(4).ColumnWidth = 42.93

# Sheet "DP": row 16 height tweak (content unchanged, purely a cosmetic resize)
$wsDp = $wb.Worksheets.Item("DP")
$wsDp.Rows.Item(16).RowHeight = 28.5
